$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The phone-number column (B2:B6) was stored as a plain number
# (989123619975, leading "0" lost) under a numeric format (numFmtId 1).
# Re-enter it as the text value "09123619975" and switch the column to a
# text number format ("@", numFmtId 49) so the leading zero is preserved.
$phone = "09123619975"
$rng = $ws.Range("B2:B6")

# Apply the text format *before* writing the values, otherwise Excel will
# auto-coerce the digit string back into a number and drop the leading 0.
$rng.NumberFormat = "@"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 2).Value = $phone
}

# The author's selection moved from A11 (stale/out of range) to B2.
$ws.Range("B2").Select()
